$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1671.0714
$ws.Range("I113").Value = 1599.4445
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 1599.4445
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 1654.5555
$ws.Range("N113").Value = -8308
$ws.Range("H116").Value = 5678.686
$ws.Range("I116").Value = 6970.174
$ws.Range("J116").Value = 4617.8213
$ws.Range("K116").Value = 6970.174
$ws.Range("L116").Value = 4617.8213
$ws.Range("M116").Value = -3528.174
$ws.Range("N116").Value = -11501.8213
$ws.Range("H129").Value = 889.0345
$ws.Range("I129").Value = 494.57144
$ws.Range("J129").Value = 1014.5455
$ws.Range("K129").Value = 1483.71432
$ws.Range("L129").Value = 3043.6365
$ws.Range("M129").Value = 3516.28568
$ws.Range("N129").Value = -13043.6365
$ws.Range("H132").Value = 254999.97
$ws.Range("I132").Value = 5096.943
$ws.Range("J132").Value = 2004321.2
$ws.Range("K132").Value = 15290.829
$ws.Range("L132").Value = 6012963.6
$ws.Range("M132").Value = -12760.829
$ws.Range("N132").Value = -6018023.6
$ws.Range("H138").Value = 106010.97
$ws.Range("I138").Value = 2545.077
$ws.Range("J138").Value = 121835.164
$ws.Range("K138").Value = 7635.231000000001
$ws.Range("L138").Value = 365505.492
$ws.Range("M138").Value = -2495.231000000001
$ws.Range("N138").Value = -375785.492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M30").ClearContents()
$ws.Range("H30").Value = 709
$ws.Range("I30").Value = 709
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 709
$ws.Range("L30").Value = 0
$ws.Range("N30").Value = -559
$ws.Range("H74").Value = 28479.27
$ws.Range("I74").Value = 30815.266
$ws.Range("J74").Value = 2004.6666
$ws.Range("K74").Value = 30815.266
$ws.Range("L74").Value = 2004.6666
$ws.Range("M74").Value = -29941.266
$ws.Range("N74").Value = -3752.6666
$ws.Range("H77").Value = 28479.27
$ws.Range("I77").Value = 30815.266
$ws.Range("J77").Value = 2004.6666
$ws.Range("K77").Value = 154076.33
$ws.Range("L77").Value = 10023.333
$ws.Range("M77").Value = -149708.33
$ws.Range("N77").Value = -18759.333
$ws.Range("H86").Value = 1000000000
$ws.Range("J86").Value = 1000000000
$ws.Range("L86").Value = 1000000000
$ws.Range("N86").Value = -1000002372
$ws.Range("H89").Value = 1000000000
$ws.Range("J89").Value = 1000000000
$ws.Range("L89").Value = 3000000000
$ws.Range("N89").Value = -3000011856
$ws.Range("H122").Value = 736.3333
$ws.Range("I122").Value = 723.15
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2169.45
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 280.5500000000002
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 80010
$ws.Range("J23").Value = 80010
$ws.Range("L23").Value = 80010
$ws.Range("N23").Value = -80490
$ws.Range("H27").Value = 80010
$ws.Range("J27").Value = 80010
$ws.Range("L27").Value = 80010
$ws.Range("N27").Value = -80394
$ws.Range("H50").Value = 9859.333000000001
$ws.Range("J50").Value = 10119.272
$ws.Range("L50").Value = 10119.272
$ws.Range("N50").Value = -11369.272
$ws.Range("H60").Value = 8035.478
$ws.Range("J60").Value = 8264.362999999999
$ws.Range("L60").Value = 8264.362999999999
$ws.Range("N60").Value = -9286.362999999999
$ws.Range("H114").Value = 100684
$ws.Range("J114").Value = 100684
$ws.Range("L114").Value = 100684
$ws.Range("N114").Value = -109362
$ws.Range("H133").Value = 36556.125
$ws.Range("J133").Value = 36556.125
$ws.Range("L133").Value = 36556.125
$ws.Range("N133").Value = -41616.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 518.1786
$ws.Range("I113").Value = 455.06668
$ws.Range("J113").Value = 591
$ws.Range("K113").Value = 1365.20004
$ws.Range("L113").Value = 1773
$ws.Range("M113").Value = 804.7999599999998
$ws.Range("N113").Value = -6113
$ws.Range("H131").Value = 157074.69
$ws.Range("J131").Value = 179460.36
$ws.Range("L131").Value = 538381.08
$ws.Range("N131").Value = -548461.08

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2098.7036
$ws.Range("I122").Value = 1951.3684
$ws.Range("J122").Value = 2448.625
$ws.Range("K122").Value = 5854.1052
$ws.Range("L122").Value = 7345.875
$ws.Range("M122").Value = -3404.1052
$ws.Range("N122").Value = -12245.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 14217
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 17000.4
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 17000.4
$ws.Range("M2").Value = -188
$ws.Range("N2").Value = -17224.4
$ws.Range("H25").Value = 5502.3335
$ws.Range("I25").Value = 753.5
$ws.Range("J25").Value = 15000
$ws.Range("K25").Value = 753.5
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = -523.5
$ws.Range("N25").Value = -15460
$ws.Range("H68").Value = 1568.92
$ws.Range("I68").Value = 1144.4445
$ws.Range("J68").Value = 2660.4285
$ws.Range("K68").Value = 1144.4445
$ws.Range("L68").Value = 2660.4285
$ws.Range("M68").Value = -395.4445000000001
$ws.Range("N68").Value = -4158.4285
$ws.Range("H71").Value = 1568.92
$ws.Range("I71").Value = 1144.4445
$ws.Range("J71").Value = 2660.4285
$ws.Range("K71").Value = 5722.2225
$ws.Range("L71").Value = 13302.1425
$ws.Range("M71").Value = -1978.2225
$ws.Range("N71").Value = -20790.1425
$ws.Range("H122").Value = 2698.1667
$ws.Range("I122").Value = 2749.318
$ws.Range("J122").Value = 2617.7856
$ws.Range("K122").Value = 8247.954000000002
$ws.Range("L122").Value = 7853.3568
$ws.Range("M122").Value = -5797.954000000002
$ws.Range("N122").Value = -12753.3568
$ws.Range("H132").Value = 200025.8
$ws.Range("I132").Value = 49250.49
$ws.Range("J132").Value = 560211.25
$ws.Range("K132").Value = 147751.47
$ws.Range("L132").Value = 1680633.75
$ws.Range("M132").Value = -145221.47
$ws.Range("N132").Value = -1685693.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1698340
$ws.Range("I136").Value = 2382955.2
$ws.Range("J136").Value = 557314.4399999999
$ws.Range("K136").Value = 7148865.600000001
$ws.Range("L136").Value = 1671943.32
$ws.Range("M136").Value = -7146315.600000001
$ws.Range("N136").Value = -1677043.32
